$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

$ws.Range("H28").Value = 312.14285
$ws.Range("I28").Value = 131.66667
$ws.Range("K28").Value = 131.66667
$ws.Range("M28").Value = 353.33333

$ws.Range("H33").Value = 673.8461
$ws.Range("I33").Value = 669.2727
$ws.Range("K33").Value = 669.2727
$ws.Range("M33").Value = -440.2727

$ws.Range("H34").Value = 54500
$ws.Range("I34").Value = 8000
$ws.Range("K34").Value = 8000
$ws.Range("M34").Value = -7797

$ws.Range("H36").Value = 54500
$ws.Range("I36").Value = 8000
$ws.Range("K36").Value = 8000
$ws.Range("M36").Value = -7285

$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("N79").ClearContents()

$ws.Range("H103").Value = 2225
$ws.Range("J103").Value = 3300
$ws.Range("L103").Value = 9900
$ws.Range("N103").Value = -11072

$ws.Range("H137").Value = 1464.1818
$ws.Range("I137").Value = 907.8570999999999
$ws.Range("K137").Value = 2723.5713
$ws.Range("M137").Value = -173.5712999999996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5270.2666
$ws.Range("I32").Value = 5270.2666
$ws.Range("K32").Value = 5270.2666
$ws.Range("M32").Value = -4983.2666

$ws.Range("H132").Value = 1904.7778
$ws.Range("I132").Value = 1751.4286
$ws.Range("J132").Value = 2441.5
$ws.Range("K132").Value = 5254.2858
$ws.Range("L132").Value = 7324.5
$ws.Range("M132").Value = -2724.2858
$ws.Range("N132").Value = -12384.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 2070.8333
$ws.Range("J5").Value = 2608.3333
$ws.Range("L5").Value = 2608.3333
$ws.Range("N5").Value = -2834.3333

$ws.Range("H10").Value = 568.3333
$ws.Range("I10").Value = 568.3333
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 568.3333
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -428.3333
$ws.Range("N10").ClearContents()

$ws.Range("H20").Value = 1127.6666
$ws.Range("I20").Value = 990
$ws.Range("K20").Value = 990
$ws.Range("M20").Value = -743

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 3000
$ws.Range("J10").Value = 3500
$ws.Range("L10").Value = 3500
$ws.Range("N10").Value = -3778

$ws.Range("H31").Value = 4809.5557
$ws.Range("J31").Value = 6875
$ws.Range("L31").Value = 6875
$ws.Range("N31").Value = -7465

$ws.Range("H34").Value = 4809.5557
$ws.Range("J34").Value = 6875
$ws.Range("L34").Value = 6875
$ws.Range("N34").Value = -7279

$ws.Range("H62").Value = 4281.2
$ws.Range("I62").Value = 3999.5
$ws.Range("J62").Value = 4469
$ws.Range("K62").Value = 3999.5
$ws.Range("L62").Value = 4469
$ws.Range("M62").Value = -3375.5
$ws.Range("N62").Value = -5717

$ws.Range("H65").Value = 4281.2
$ws.Range("I65").Value = 3999.5
$ws.Range("J65").Value = 4469
$ws.Range("K65").Value = 19997.5
$ws.Range("L65").Value = 22345
$ws.Range("M65").Value = -16877.5
$ws.Range("N65").Value = -28585

$ws.Range("H103").Value = 38305.832
$ws.Range("I103").Value = 38305.832
$ws.Range("K103").Value = 38305.832
$ws.Range("M103").Value = -37133.832

$ws.Range("H105").Value = 1999
$ws.Range("I105").Value = 1998
$ws.Range("K105").Value = 1998
$ws.Range("M105").Value = -251

$ws.Range("H132").Value = 2188.9
$ws.Range("I132").Value = 2048.625
$ws.Range("J132").Value = 2750
$ws.Range("K132").Value = 6145.875
$ws.Range("L132").Value = 8250
$ws.Range("M132").Value = -3615.875
$ws.Range("N132").Value = -13310

$ws.Range("H134").Value = 2399.6
$ws.Range("I134").Value = 1999.5
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 5998.5
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -3463.5
$ws.Range("N134").Value = -17070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 311.25
$ws.Range("I17").Value = 48
$ws.Range("J17").Value = 750
$ws.Range("K17").Value = 144
$ws.Range("L17").Value = 2250
$ws.Range("M17").Value = 25
$ws.Range("N17").Value = -2588

$ws.Range("H39").Value = 1674.6666
$ws.Range("J39").Value = 1959.6
$ws.Range("L39").Value = 5878.799999999999
$ws.Range("N39").Value = -6466.799999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 250
$ws.Range("I29").Value = 250
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 250
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 40
$ws.Range("N29").ClearContents()

$ws.Range("H132").Value = 1000
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 3000
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -8060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 7380
$ws.Range("I20").Value = 7380
$ws.Range("K20").Value = 7380
$ws.Range("M20").Value = -7154

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 226.75
$ws.Range("I100").Value = 187.71428
$ws.Range("K100").Value = 375.42856
$ws.Range("M100").Value = 165.57144

$ws.Range("H113").Value = 666
$ws.Range("J113").Value = 499
$ws.Range("L113").Value = 1497
$ws.Range("N113").Value = -5837

$ws.Range("H122").Value = 1650
$ws.Range("I122").Value = 1650
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4950
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2500
$ws.Range("N122").ClearContents()
